# Update column headers: PREREQ_COURSES -> Prerequisites, COREQ_COURSES -> CoRequisites
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E1").Value = "Prerequisites"
$ws.Range("F1").Value = "CoRequisites"

# Restore selection to the header row only (A1:I1) as left by the author after editing
$ws.Range("A1:I1").Select()
